# Auto-generated: apply numeric corrections to the 'profits' recompute columns
# (H/I/J/K/L/M/N) across all 8 item-type sheets, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$sheetEdits = @{}
$sheetEdits["ALC"] = @{
    "H18" = 5833.9
    "I18" = 5833.9
    "K18" = 5833.9
    "M18" = -5549.9
    "H33" = 347.6154
    "I33" = 259.375
    "K33" = 259.375
    "M33" = -30.375
    "H40" = 2448.3076
    "I40" = 2440.8572
    "J40" = 2457
    "K40" = 2440.8572
    "L40" = 2457
    "M40" = -2265.8572
    "N40" = -2807
    "H74" = 5457.9585
    "I74" = 5142.4287
    "K74" = 5142.4287
    "M74" = -4206.4287
    "H77" = 5457.9585
    "I77" = 5142.4287
    "K77" = 25712.1435
    "M77" = -21032.1435
    "H86" = 6494.857
    "I86" = 6693
    "K86" = 6693
    "M86" = -5570
    "H88" = 2517.2942
    "I88" = 2033.3334
    "J88" = 3061.75
    "K88" = 2033.3334
    "L88" = 3061.75
    "M88" = -1627.3334
    "N88" = -3873.75
    "H89" = 6494.857
    "I89" = 6693
    "K89" = 33465
    "M89" = -27849
    "H91" = 2517.2942
    "I91" = 2033.3334
    "J91" = 3061.75
    "K91" = 2033.3334
    "L91" = 3061.75
    "M91" = -629.3334
    "N91" = -5869.75
    "H97" = 1492.6666
    "J97" = 1492.6666
    "L97" = 4477.9998
    "N97" = -5469.9998
    "H100" = 2599.4
    "I100" = 1749.5
    "J100" = 3166
    "K100" = 1749.5
    "L100" = 3166
    "M100" = -1208.5
    "N100" = -4248
    "H112" = 3998.5
    "J112" = 3998.5
    "L112" = 11995.5
    "N112" = -14211.5
    "H141" = 2740.1667
    "I141" = 2288.4666
    "J141" = 4998.6665
    "K141" = 6865.399800000001
    "L141" = 14995.9995
    "M141" = -1685.399800000001
    "N141" = -25355.9995
}

$sheetEdits["ARM"] = @{
    "H32" = 192671.03
    "I32" = 220273.16
    "J32" = 11285.714
    "K32" = 220273.16
    "L32" = 11285.714
    "M32" = -219986.16
    "N32" = -11859.714
    "H45" = 3528.2856
    "I45" = 3400
    "J45" = 3549.6667
    "K45" = 3400
    "L45" = 3549.6667
    "M45" = -3023
    "N45" = -4303.6667
    "H97" = 38462940
    "I97" = 1314.7273
    "J97" = 250001890
    "K97" = 1314.7273
    "L97" = 250001890
    "M97" = -818.7273
    "N97" = -250002882
    "H132" = 2086400.6
    "I132" = 2275619
    "K132" = 6826857
    "M132" = -6824327
}

$sheetEdits["BSM"] = @{
    "H28" = 46500
    "J28" = 46500
    "L28" = 46500
    "N28" = -47088
    "H107" = 1247.8182
    "I107" = 1247.8182
    "K107" = 1247.8182
    "M107" = 672.1818000000001
    "H132" = 98035
    "J132" = 98035
    "L132" = 98035
    "N132" = -108155
    "H134" = 3252.2856
    "I134" = 3470.2307
    "J134" = 2898.125
    "K134" = 10410.6921
    "L134" = 8694.375
    "M134" = -7875.6921
    "N134" = -13764.375
}

$sheetEdits["CRP"] = @{
    "H58" = 5356.863
    "I58" = 4284.0454
    "J58" = 6170.724
    "K58" = 4284.0454
    "L58" = 6170.724
    "M58" = -4081.0454
    "N58" = -6576.724
    "J62" = 0
    "L62" = 0
    "N62" = $null
    "J65" = 0
    "L65" = 0
    "N65" = $null
    "H105" = 2222
    "I105" = 571.2857
    "J105" = 7999.5
    "K105" = 571.2857
    "L105" = 7999.5
    "M105" = 1175.7143
    "N105" = -11493.5
    "H134" = 1699.5294
    "I134" = 1251.0968
    "K134" = 3753.2904
    "M134" = -1218.2904
    "H136" = 5356.863
    "I136" = 4284.0454
    "J136" = 6170.724
    "K136" = 12852.1362
    "L136" = 18512.172
    "M136" = -10302.1362
    "N136" = -23612.172
}

$sheetEdits["CUL"] = @{
    "H7" = 16.75
    "I7" = 19.333334
    "J7" = 9
    "K7" = 58.000002
    "L7" = 27
    "M7" = 53.999998
    "N7" = -251
    "H9" = 533501
    "J9" = 600002
    "L9" = 1800006
    "N9" = -1800454
}

$sheetEdits["GSM"] = @{
    "H80" = 2196.5454
    "I80" = 1794.5
    "J80" = 3268.6667
    "K80" = 1794.5
    "L80" = 3268.6667
    "M80" = -796.5
    "N80" = -5264.6667
    "H83" = 2196.5454
    "I83" = 1794.5
    "J83" = 3268.6667
    "K83" = 8972.5
    "L83" = 16343.3335
    "M83" = -3980.5
    "N83" = -26327.3335
}

$sheetEdits["LTW"] = @{
    "H61" = 14988
    "I61" = 17605
    "K61" = 17605
    "M61" = -17403
    "H82" = 1461
    "I82" = 1691.5
    "K82" = 1691.5
    "M82" = -1330.5
    "H85" = 1461
    "I85" = 1691.5
    "K85" = 1691.5
    "M85" = -443.5
    "H113" = 14988
    "I113" = 17605
    "K113" = 17605
    "M113" = -15435
    "H122" = 9989.933999999999
    "I122" = 12028.223
    "J122" = 6932.5
    "K122" = 36084.669
    "L122" = 20797.5
    "M122" = -33634.669
    "N122" = -25697.5
}

$sheetEdits["WVR"] = @{
    "H81" = 3158.1538
    "J81" = 5945.6665
    "L81" = 11891.333
    "N81" = -14013.333
    "H84" = 3158.1538
    "J84" = 5945.6665
    "L84" = 59456.665
    "N84" = -70064.66500000001
    "H96" = 3308.8
    "I96" = 3102.9092
    "K96" = 3102.9092
    "M96" = -1729.9092
    "H107" = 1826.2307
    "I107" = 1074.2273
    "K107" = 3222.6819
    "M107" = -1302.6819
}

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $edits = $sheetEdits[$sheetName]
    foreach ($cellRef in $edits.Keys) {
        $ws.Range($cellRef).Value = $edits[$cellRef]
    }
}

Write-Host "Applied $($sheetEdits.Values | ForEach-Object { $_.Keys.Count } | Measure-Object -Sum | Select-Object -ExpandProperty Sum) cell edits across $($sheetEdits.Keys.Count) sheets."
